# Add Team Chat functionality with XLSX storage and @mentions support
#
# Adds a new "Chat" worksheet at the end of the workbook (after "Vendors")
# containing a log of team chat messages, including @mention style
# recipients. Timestamps are compact numeric strings (yyyyMMddHHmmss) that
# must be preserved as TEXT (not coerced to a Number) - we force that with
# the classic leading-apostrophe "store as text" trick, same as Excel's UI.

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately after the current last sheet so it lands
# at the very end of the tab strip (Worksheets.Add with no "Before" puts it
# first by default).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Chat"

# Header row
$headers = @("Timestamp", "User", "Message", "Type", "Recipients", "Status")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Chat log rows: Timestamp, User, Message, Type, Recipients, Status
$rows = @(
    @("20241201143000", "Alyssa",   "Hey team, where are we on the Johnson case?",          "question",  "all",       "active"),
    @("20241201143100", "Dr. Moore","I just reviewed the medication list, all looks good",  "update",    "all",       "active"),
    @("20241201143200", "Christa",  "Family meeting scheduled for tomorrow at 2pm",          "info",      "all",       "active"),
    @("20241201143300", "Amber",    "Insurance approval came through!",                      "good_news", "all",       "active"),
    @("20241201143400", "Alyssa",   "Great work everyone!",                                  "comment",   "all",       "active"),
    @("20241201143500", "Dr. Moore","@Christa - can you prep the meeting notes?",            "task",      "Christa",   "active"),
    @("20241201143600", "Christa",  "On it! Will have them ready by EOD",                    "response",  "Dr. Moore", "active"),
    @("20241201143700", "Amber",    "@Alyssa - need your input on the billing codes",        "question",  "Alyssa",    "active"),
    @("20241201143800", "Alyssa",   "I'll review and get back to you by 5pm",                "response",  "Amber",     "active")
)

$rowIndex = 2
foreach ($row in $rows) {
    # Timestamp: prefix with an apostrophe so it is stored as text, not a
    # number - mirrors the "Numbers Stored as Text" condition in the sheet.
    $ws.Cells.Item($rowIndex, 1).Value = "'" + $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
